$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (table shrinks from 18 to 17 data rows)
$ws.Rows.Item(19).Delete()

# Overwrite the remaining data rows (2-18) with the new roster data
$ws.Range("A2").Value = "Cade Cunningham"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Detroit Pistons"
$ws.Range("A3").Value = "Coby White"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Chicago Bulls"
$ws.Range("A4").Value = "Ausar Thompson"
$ws.Range("B4").Value = "SF,PF"
$ws.Range("C4").Value = "Detroit Pistons"
$ws.Range("A5").Value = "Derrick White"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Boston Celtics"
$ws.Range("A6").Value = "Carlton Carrington"
$ws.Range("B6").Value = "PG,SG"
$ws.Range("C6").Value = "Washington Wizards"
$ws.Range("A7").Value = "LaMelo Ball"
$ws.Range("B7").Value = "PG,SG"
$ws.Range("C7").Value = "Charlotte Hornets"
$ws.Range("A8").Value = "Malcolm Brogdon"
$ws.Range("B8").Value = "PG,SG"
$ws.Range("C8").Value = "Washington Wizards"
$ws.Range("A9").Value = "Devin Vassell"
$ws.Range("B9").Value = "SG,SF"
$ws.Range("C9").Value = "San Antonio Spurs"
$ws.Range("A10").Value = "Naz Reid"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Minnesota Timberwolves"
$ws.Range("A11").Value = "Onyeka Okongwu"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Atlanta Hawks"
$ws.Range("A12").Value = "Malik Monk"
$ws.Range("B12").Value = "PG,SG,SF"
$ws.Range("C12").Value = "Sacramento Kings"
$ws.Range("A13").Value = "Isaiah Hartenstein"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Oklahoma City Thunder"
$ws.Range("A14").Value = "Damian Lillard"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Milwaukee Bucks"
$ws.Range("A15").Value = "Jaxson Hayes"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Los Angeles Lakers"
$ws.Range("A16").Value = "Collin Sexton"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Utah Jazz"
$ws.Range("A17").Value = "Julius Randle"
$ws.Range("B17").Value = "PF,C"
$ws.Range("C17").Value = "Minnesota Timberwolves"
$ws.Range("A18").Value = "Anthony Davis"
$ws.Range("B18").Value = "PF,C"
$ws.Range("C18").Value = "Dallas Mavericks"
